$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 12 (2536_BLB11ACC_r1.fq) - entire row removed
$ws.Rows.Item(12).Delete()

# 2. Fix B22 (2536_BLB11MCC_r1.fq) value - corrected read count
$ws.Range("B22").Value = 39956649

# 3. Add header cells for new columns (string order matches the commit's sharedStrings layout)
$ws.Range("B1").Value = "raw data fastqc"
$ws.Range("C1").Value = "trim galore fastqc"

# 4. Add new row 33 for 2536_BLB11VMH_r1.fq
$ws.Range("A33").Value = "2536_BLB11VMH_r1.fq"
$ws.Range("B33").Value = 24638802

# 5. Fill column C (trim galore fastqc counts) for all data rows
$ws.Range("C2").Value = 32199711
$ws.Range("C3").Value = 44235205
$ws.Range("C4").Value = 35656479
$ws.Range("C5").Value = 34045832
$ws.Range("C6").Value = 26738571
$ws.Range("C7").Value = 45433468
$ws.Range("C8").Value = 15365036
$ws.Range("C9").Value = 23535898
$ws.Range("C10").Value = 12782594
$ws.Range("C11").Value = 24655684
$ws.Range("C12").Value = 26869820
$ws.Range("C13").Value = 16345766
$ws.Range("C14").Value = 25387628
$ws.Range("C15").Value = 22213151
$ws.Range("C16").Value = 22465755
$ws.Range("C17").Value = 36939954
$ws.Range("C18").Value = 19039160
$ws.Range("C19").Value = 29892525
$ws.Range("C20").Value = 21620848
$ws.Range("C21").Value = 24766166
$ws.Range("C22").Value = 36359217
$ws.Range("C23").Value = 15269064
$ws.Range("C24").Value = 21414984
$ws.Range("C25").Value = 28936004
$ws.Range("C26").Value = 22598492
$ws.Range("C27").Value = 19524250
$ws.Range("C28").Value = 36065618
$ws.Range("C29").Value = 34039803
$ws.Range("C30").Value = 17638909
$ws.Range("C31").Value = 25835203
$ws.Range("C32").Value = 11129988
$ws.Range("C33").Value = 23835493

# 6. Apply explicit black font color to C18 (matches author formatting)
$ws.Range("C18").Font.Color = 0

# 7. Column widths for the new columns
$ws.Columns.Item(2).ColumnWidth = 16.5
$ws.Columns.Item(3).ColumnWidth = 17.5

# 8. View state cosmetics (zoom + selection)
$excel.ActiveWindow.Zoom = 114
$ws.Range("C36").Select()
